# The workbook contains a daily price log for "Zapallo italiano" at
# Femacal de La Calera. A new daily record was inserted as row 384
# (pushing the former rows 384-473 down to 385-474), growing the sheet
# from A1:R473 to A1:R474.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 384; this shifts rows 384:473
# down to 385:474 and extends the used range accordingly.
$ws.Rows("384").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(384, 1).Value  = 3
$ws.Cells.Item(384, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(384, 3).Value  = "Coquimbo"
$ws.Cells.Item(384, 4).Value  = 44711
$ws.Cells.Item(384, 5).Value  = 5
$ws.Cells.Item(384, 6).Value  = 100112032
$ws.Cells.Item(384, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(384, 8).Value  = "Sin especificar"
$ws.Cells.Item(384, 9).Value  = "Primera"
$ws.Cells.Item(384, 10).Value = 205
$ws.Cells.Item(384, 11).Value = 11500
$ws.Cells.Item(384, 12).Value = 12000
$ws.Cells.Item(384, 13).Value = 11671
$ws.Cells.Item(384, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(384, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(384, 16).Value = 195
$ws.Cells.Item(384, 17).Value = 60
$ws.Cells.Item(384, 18).Value = "Hortaliza"
